$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (JOBY) updates
$ws.Range("K2").Value = 57.8
$ws.Range("N2").Value = 54.82400714602223

# Row 3 (ACHR) updates
$ws.Range("K3").Value = 54.4
$ws.Range("N3").Value = 54.82400714602223
